{"js": "// The document body is a single empty paragraph. The edit sets the\n// paragraph mark's run language to English (US), which Word persists as\n//   <w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>\n// on that (otherwise empty) paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\nconst range = paragraph.getRange();\nrange.languageId = \"en-US\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document body is a single empty paragraph. The edit sets the\n# paragraph mark's run language to English (US), which Word stores as\n# <w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr> on that paragraph.\n$p = $d.Paragraphs.Item(1)\n$p.Range.LanguageID = \"en-US\"\n"}
